$d = $word.ActiveDocument

# 1. Title year 2022 -> 2024 (second run '2' -> '4')
$d.Content.Find.Execute("2022年暑期课程", $true, $false, $false, $false, $false, $true, 1, $false, "2024年暑期课程", 2)

# 2. Course homepage URL
$d.Content.Find.Execute("http://staff.ustc.edu.cn/~renjiec/SummerSchool_2022/index.html", $true, $false, $false, $false, $false, $true, 1, $false, "https://ustc-gcl-f.github.io/course/SummerSchool_2024/index.html", 2)

# 3. Submission date year 2022 -> 2024, day 14 -> 18
$d.Content.Find.Execute("2022年8月14日", $true, $false, $false, $false, $false, $true, 1, $false, "2024年8月18日", 2)

# 4. Name and email
$d.Content.Find.Execute("曹合智同学（caohezhi21@mail.ustc.edu.cn）", $true, $false, $false, $false, $false, $true, 1, $false, "张老师（434484980@qq.com）", 2)

